# Monitoreo de actividades del 20 al 27 de mayo
# Adds two closed non-conformity rows (27 and 28) to the "No Conformidades"
# report sheet, covering work monitored between 2016-05-20 and 2016-05-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30 (item #27) was a placeholder row (only "A30" had the sequence
# number). Fill in the rest of the row with the closed non-conformity that
# was tracked for this monitoring cycle. Re-using the cells that already
# exist here keeps their existing number/alignment/fill/border formatting
# untouched.
$ws.Range("B30").Value = "No todas las tareas fueron realizadas"
$ws.Range("C30").Value = "Ventas"
$ws.Range("D30").Value = 42517
$ws.Range("E30").Value = 42517
$ws.Range("F30").Value = "Cerrada"
$ws.Range("G30").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."

# --- Row 31 (item #28) is a brand-new row. Clone the formatting of row 30
# (which already carries the correct borders/fills/alignment/number
# formats for every column in this table) down into row 31, then overwrite
# the values for the new entry.
$ws.Range("A30:G30").Copy($ws.Range("A31:G31"))
$ws.Range("A31").Value = 28
$ws.Range("B31").Value = "Las tareas no fueron realizadas"
$ws.Range("C31").Value = "Compras"
$ws.Range("D31").Value = 42517
$ws.Range("E31").Value = 42517
$ws.Range("F31").Value = "Cerrada"
$ws.Range("G31").Value = "Esta anomalia solo se reporta ya que son tareas ya realizadas. Se presentan para futuras actividades y se escalan."

# Both new rows hold long wrapped comments in column G, so they need the
# same taller row height used by the other commented rows in the sheet.
$ws.Rows("30:31").RowHeight = 75

# Extend the STATUS column's list data-validation so it also covers the
# freshly added row 31 (it previously stopped at row 30).
$ws.Range("F4:F31").Validation.Delete()
$ws.Range("F4:F31").Validation.Add(3, 1, 1, '"En proceso,Cerrada,Cancelada,Rechazada"', "0")

# Scroll/select to where the new rows are, matching where the author left
# the cursor after the edit.
$excel.ActiveWindow.ScrollRow = 29
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C34").Select()
